# The deck originally ships two themes:
#   theme1.xml -> "Integral"      (used by the slide master / all slides)
#   theme2.xml -> "Office Theme"  (used by the notes master)
#
# The target edit swaps the two themes' colour schemes so that the slide
# master (the one driving the visible deck) now carries the stock
# "Office Theme" palette instead of "Integral".
#
# Office Theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink):
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
